$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1592.3077
$ws.Range("I40").Value = 1550
$ws.Range("J40").Value = 1660
$ws.Range("K40").Value = 1550
$ws.Range("L40").Value = 1660
$ws.Range("M40").Value = -1375
$ws.Range("N40").Value = -2010

# Row 43
$ws.Range("H43").Value = 1633.1666
$ws.Range("I43").Value = 950
$ws.Range("J43").Value = 1974.75
$ws.Range("K43").Value = 950
$ws.Range("L43").Value = 1974.75
$ws.Range("M43").Value = -881
$ws.Range("N43").Value = -2112.75

# Row 106
$ws.Range("H106").Value = 541.2727
$ws.Range("I106").Value = 395.4
$ws.Range("K106").Value = 395.4
$ws.Range("M106").Value = 235.6

# Row 138
$ws.Range("H138").Value = 2442643.8
$ws.Range("I138").Value = 14287724
$ws.Range("J138").Value = 3950.647
$ws.Range("K138").Value = 42863172
$ws.Range("L138").Value = 11851.941
$ws.Range("M138").Value = -42858032
$ws.Range("N138").Value = -22131.941

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 18007.568
$ws.Range("I32").Value = 18959.342
$ws.Range("K32").Value = 18959.342
$ws.Range("M32").Value = -18672.342

# Row 61
$ws.Range("H61").Value = 2277.75
$ws.Range("I61").Value = 1903.6666
$ws.Range("J61").Value = 3400
$ws.Range("K61").Value = 1903.6666
$ws.Range("L61").Value = 3400
$ws.Range("M61").Value = -1691.6666
$ws.Range("N61").Value = -3824

# Row 74
$ws.Range("H74").Value = 1269.6818
$ws.Range("I74").Value = 1124.0555
$ws.Range("J74").Value = 1925
$ws.Range("K74").Value = 1124.0555
$ws.Range("L74").Value = 1925
$ws.Range("M74").Value = -250.0554999999999
$ws.Range("N74").Value = -3673

# Row 77
$ws.Range("H77").Value = 1269.6818
$ws.Range("I77").Value = 1124.0555
$ws.Range("J77").Value = 1925
$ws.Range("K77").Value = 5620.2775
$ws.Range("L77").Value = 9625
$ws.Range("M77").Value = -1252.2775
$ws.Range("N77").Value = -18361

# Row 92
$ws.Range("H92").Value = 199989.67
$ws.Range("J92").Value = 199989.67
$ws.Range("L92").Value = 199989.67
$ws.Range("N92").Value = -204981.67

# Row 136
$ws.Range("H136").Value = 2277.75
$ws.Range("I136").Value = 1903.6666
$ws.Range("J136").Value = 3400
$ws.Range("K136").Value = 5710.9998
$ws.Range("L136").Value = 10200
$ws.Range("M136").Value = -3160.9998
$ws.Range("N136").Value = -15300

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 88
$ws.Range("H88").Value = 30949.75
$ws.Range("J88").Value = 30949.75
$ws.Range("L88").Value = 30949.75
$ws.Range("N88").Value = -31761.75

# Row 91
$ws.Range("H91").Value = 30949.75
$ws.Range("J91").Value = 30949.75
$ws.Range("L91").Value = 30949.75
$ws.Range("N91").Value = -33757.75

# Row 92
$ws.Range("H92").Value = 152798
$ws.Range("J92").Value = 152798
$ws.Range("L92").Value = 152798
$ws.Range("N92").Value = -157790

# Row 95
$ws.Range("H95").Value = 40089.145
$ws.Range("J95").Value = 40089.145
$ws.Range("L95").Value = 40089.145
$ws.Range("N95").Value = -45581.145

# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# Row 107
$ws.Range("H107").Value = 35777.734
$ws.Range("J107").Value = 1580
$ws.Range("L107").Value = 1580
$ws.Range("N107").Value = -5420

# Row 139
$ws.Range("H139").Value = 105677.14
$ws.Range("J139").Value = 105677.14
$ws.Range("L139").Value = 105677.14
$ws.Range("N139").Value = -115957.14

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 20836618
$ws.Range("I31").Value = 43481080
$ws.Range("J31").Value = 3715.2
$ws.Range("K31").Value = 43481080
$ws.Range("L31").Value = 3715.2
$ws.Range("M31").Value = -43480785
$ws.Range("N31").Value = -4305.2

# Row 34
$ws.Range("H34").Value = 20836618
$ws.Range("I34").Value = 43481080
$ws.Range("J34").Value = 3715.2
$ws.Range("K34").Value = 43481080
$ws.Range("L34").Value = 3715.2
$ws.Range("M34").Value = -43480878
$ws.Range("N34").Value = -4119.2

# Row 88
$ws.Range("H88").Value = 32298.143
$ws.Range("J88").Value = 32298.143
$ws.Range("L88").Value = 32298.143
$ws.Range("N88").Value = -33110.143

# Row 91
$ws.Range("H91").Value = 32298.143
$ws.Range("J91").Value = 32298.143
$ws.Range("L91").Value = 32298.143
$ws.Range("N91").Value = -35106.143

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 297.35
$ws.Range("I17").Value = 289.6
$ws.Range("J17").Value = 299.93332
$ws.Range("K17").Value = 868.8000000000001
$ws.Range("L17").Value = 899.7999599999999
$ws.Range("M17").Value = -699.8000000000001
$ws.Range("N17").Value = -1237.79996

# Row 104
$ws.Range("H104").Value = 7000
$ws.Range("J104").Value = 7000
$ws.Range("L104").Value = 21000
$ws.Range("N104").Value = -26242

# Row 131
$ws.Range("H131").Value = 45461116
$ws.Range("J131").Value = 76926260
$ws.Range("L131").Value = 230778780
$ws.Range("N131").Value = -230788860

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 64
$ws.Range("H64").Value = 29999.285
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 29999.285
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 29999.285
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -30495.285

# Row 67
$ws.Range("H67").Value = 29999.285
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 29999.285
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 29999.285
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -31715.285

# Row 102
$ws.Range("H102").Value = 2148.76
$ws.Range("I102").Value = 1890.7368
$ws.Range("K102").Value = 1890.7368
$ws.Range("M102").Value = -268.7367999999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 150.25
$ws.Range("I55").Value = 150.25
$ws.Range("K55").Value = 150.25
$ws.Range("M55").Value = 22.75

# Row 97
$ws.Range("H97").Value = 25845.615
$ws.Range("J97").Value = 25845.615
$ws.Range("L97").Value = 25845.615
$ws.Range("N97").Value = -27827.615

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

# Row 92
$ws.Range("H92").Value = 25183.334
$ws.Range("J92").Value = 25183.334
$ws.Range("L92").Value = 25183.334
$ws.Range("N92").Value = -30175.334

# Row 95
$ws.Range("H95").Value = 101562.5
$ws.Range("J95").Value = 101562.5
$ws.Range("L95").Value = 101562.5
$ws.Range("N95").Value = -107054.5
